$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 30: designator list in C30 loses "R58" (component removed from BOM)
# ---------------------------------------------------------------------------
$cell30 = $ws.Range("C30")
$cell30.Value = "R10,R13,R16,R19,R21,R23,R24,R29,R30,R39,R40,R50,R51,R57,R59,R60,R62,R64"

$green = 5287936   # RGB(0,176,80) in BGR COM ordering
$black = 0         # RGB(0,0,0)
$red   = 255       # RGB(255,0,0) in BGR COM ordering

function Set-RunFont($cell, $start, $len, $color) {
    $rng = $cell.Characters($start, $len)
    $rng.Font.Size = 10
    $rng.Font.Name = "Liberation Sans"
    $rng.Font.Color = $color
}

# Re-apply the original rich-text run formatting that existed around the
# designators that stay colored (R39, R40, R59, R60 green; R64 red), while
# the plain designators (including the run that used to end in ",R58,")
# keep the default black text.
Set-RunFont $cell30 37 3 $green    # R39
Set-RunFont $cell30 40 1 $black    # ,
Set-RunFont $cell30 41 3 $green    # R40
Set-RunFont $cell30 44 13 $black   # ,R50,R51,R57,
Set-RunFont $cell30 57 3 $green    # R59
Set-RunFont $cell30 60 1 $black    # ,
Set-RunFont $cell30 61 3 $green    # R60
Set-RunFont $cell30 64 5 $black    # ,R62,
Set-RunFont $cell30 69 3 $red      # R64

# ---------------------------------------------------------------------------
# Row 35: designator list in C35 loses "R56" (component removed from BOM)
# ---------------------------------------------------------------------------
$ws.Range("C35").Value = "R11,R14,R17,R20,R35,R36,R37,R38,R48,R49,R55"

# Leave the selection on the cell that was last edited, matching the
# author's final cursor position in the saved workbook.
$ws.Activate() | Out-Null
$ws.Range("C35").Select() | Out-Null
